# Auto commit at 2025-08-16  8:47:09.88
#
# Updates the monthly metric snapshot on "Metrics" (new month-to-date
# figures land in B3:B6 of "today", which then roll into the running
# totals on both sheets), refreshes the "TODAY()-1" reporting date on
# "today", and moves the active-tab/selection from "today" back to
# "Metrics".

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

# ---------------------------------------------------------------------
# "today" sheet: new month-to-date inputs (B3:B6) feed the formulas
# below, so write these first.
# ---------------------------------------------------------------------
$today.Range("B3").Value = 15874.69
$today.Range("B4").Value = 13649.26
$today.Range("B5").Value = 4952.7299999999996
$today.Range("B6").Value = 619

# Running totals in B11:B22 become formulas that add the above
# month-to-date figures onto the previous cumulative totals.
$today.Range("B11").Formula = "=236440.97+B3"
$today.Range("B12").Formula = "=203725.78+B4"
$today.Range("B13").Formula = "=74560.05+B5"
$today.Range("B14").Formula = "=9416+B6"
$today.Range("B15").Formula = "=3632069.54+B3"
$today.Range("B16").Formula = "=3082440.44+B4"
$today.Range("B17").Formula = "=1038202.6+B51"
$today.Range("B18").Formula = "=140104+B6"
$today.Range("B19").Formula = "=32097393.34+B3"
$today.Range("B20").Formula = "=19112310.51+B4"
$today.Range("B21").Formula = "=11319911.5+B5"
$today.Range("B22").Formula = "=1237731+B6"

# ---------------------------------------------------------------------
# "Metrics" sheet: mirrors the same refreshed cumulative totals.
# ---------------------------------------------------------------------
$metrics.Range("B2").Value  = 252315.66
$metrics.Range("B3").Value  = 217375.04
$metrics.Range("B4").Value  = 79512.78
$metrics.Range("B5").Value  = 10035
$metrics.Range("B6").Value  = 3647944.23
$metrics.Range("B7").Value  = 3096089.6999999997
$metrics.Range("B8").Value  = 1038202.6
$metrics.Range("B9").Value  = 140723
$metrics.Range("B10").Value = 32113268.030000001
$metrics.Range("B11").Value = 19125959.770000003
$metrics.Range("B12").Value = 11324864.23
$metrics.Range("B13").Value = 1238350

# ---------------------------------------------------------------------
# Reporting date bump on "today" (A1 = TODAY()-1); re-applying the
# formula lets the engine recompute the cached value for "today".
# ---------------------------------------------------------------------
$today.Range("A1").Formula = "=TODAY()-1"

# ---------------------------------------------------------------------
# Move the active sheet / selection from "today" back to "Metrics".
# ---------------------------------------------------------------------
[void]$today.Range("B11:B25").Select()

[void]$metrics.Activate()
[void]$metrics.Range("E13").Select()
